$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 422245.75
$ws.Range("J17").Value = 438865.12
$ws.Range("L17").Value = 1316595.36
$ws.Range("N17").Value = -1316931.36
$ws.Range("H103").Value = 687.4516
$ws.Range("I103").Value = 327.44446
$ws.Range("J103").Value = 834.7273
$ws.Range("K103").Value = 982.33338
$ws.Range("L103").Value = 2504.1819
$ws.Range("M103").Value = -396.33338
$ws.Range("N103").Value = -3676.1819
$ws.Range("H116").Value = 38964616
$ws.Range("I116").Value = 50202380
$ws.Range("K116").Value = 50202380
$ws.Range("M116").Value = -50198938
$ws.Range("H138").Value = 3077.9482
$ws.Range("I138").Value = 1264.8948
$ws.Range("J138").Value = 3961.2307
$ws.Range("K138").Value = 3794.6844
$ws.Range("L138").Value = 11883.6921
$ws.Range("M138").Value = 1345.3156
$ws.Range("N138").Value = -22163.6921
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1038.5883
$ws.Range("I2").Value = 1070.4
$ws.Range("J2").Value = 800
$ws.Range("K2").Value = 1070.4
$ws.Range("L2").Value = 800
$ws.Range("M2").Value = -957.4000000000001
$ws.Range("N2").Value = -1026
$ws.Range("H74").Value = 2996.8262
$ws.Range("I74").Value = 2075.25
$ws.Range("J74").Value = 5103.2856
$ws.Range("K74").Value = 2075.25
$ws.Range("L74").Value = 5103.2856
$ws.Range("M74").Value = -1201.25
$ws.Range("N74").Value = -6851.2856
$ws.Range("H77").Value = 2996.8262
$ws.Range("I77").Value = 2075.25
$ws.Range("J77").Value = 5103.2856
$ws.Range("K77").Value = 10376.25
$ws.Range("L77").Value = 25516.428
$ws.Range("M77").Value = -6008.25
$ws.Range("N77").Value = -34252.428
$ws.Range("H97").Value = 1061.7667
$ws.Range("I97").Value = 789.94446
$ws.Range("J97").Value = 1469.5
$ws.Range("K97").Value = 789.94446
$ws.Range("L97").Value = 1469.5
$ws.Range("M97").Value = -293.94446
$ws.Range("N97").Value = -2461.5
$ws.Range("H110").Value = 57812.832
$ws.Range("I110").Value = 37902.285
$ws.Range("K110").Value = 37902.285
$ws.Range("M110").Value = -35857.285
$ws.Range("H116").Value = 1038.5883
$ws.Range("I116").Value = 1070.4
$ws.Range("J116").Value = 800
$ws.Range("K116").Value = 1070.4
$ws.Range("L116").Value = 800
$ws.Range("M116").Value = 1223.6
$ws.Range("N116").Value = -5388
$ws.Range("H132").Value = 55558372
$ws.Range("I132").Value = 83335790
$ws.Range("K132").Value = 250007370
$ws.Range("M132").Value = -250004840
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1038.5883
$ws.Range("I3").Value = 1070.4
$ws.Range("J3").Value = 800
$ws.Range("K3").Value = 1070.4
$ws.Range("L3").Value = 800
$ws.Range("M3").Value = -956.4000000000001
$ws.Range("N3").Value = -1028
$ws.Range("H107").Value = 29158.945
$ws.Range("I107").Value = 21817.285
$ws.Range("J107").Value = 54854.75
$ws.Range("K107").Value = 21817.285
$ws.Range("L107").Value = 54854.75
$ws.Range("M107").Value = -19897.285
$ws.Range("N107").Value = -58694.75
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1008.6667
$ws.Range("I16").Value = 1036.125
$ws.Range("J16").Value = 789
$ws.Range("K16").Value = 1036.125
$ws.Range("L16").Value = 789
$ws.Range("M16").Value = -749.125
$ws.Range("N16").Value = -1363
$ws.Range("H113").Value = 1008.6667
$ws.Range("I113").Value = 1036.125
$ws.Range("J113").Value = 789
$ws.Range("K113").Value = 1036.125
$ws.Range("L113").Value = 789
$ws.Range("M113").Value = 1133.875
$ws.Range("N113").Value = -5129
$ws.Range("H122").Value = 2658.2666
$ws.Range("J122").Value = 3110.2856
$ws.Range("L122").Value = 9330.856800000001
$ws.Range("N122").Value = -14230.8568
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 18320698
$ws.Range("I4").Value = 35804820
$ws.Range("J4").Value = 4000
$ws.Range("K4").Value = 107414460
$ws.Range("L4").Value = 12000
$ws.Range("M4").Value = -107414348
$ws.Range("N4").Value = -12224
$ws.Range("H8").Value = 4333.3335
$ws.Range("I8").Value = 4333.3335
$ws.Range("K8").Value = 13000.0005
$ws.Range("M8").Value = -12861.0005
$ws.Range("H23").Value = 46.857143
$ws.Range("I23").Value = 11
$ws.Range("J23").Value = 52.833332
$ws.Range("K23").Value = 33
$ws.Range("L23").Value = 158.499996
$ws.Range("M23").Value = 202
$ws.Range("N23").Value = -628.499996
$ws.Range("H38").Value = 112.72727
$ws.Range("I38").Value = 131.11111
$ws.Range("J38").Value = 30
$ws.Range("K38").Value = 393.33333
$ws.Range("L38").Value = 90
$ws.Range("M38").Value = -46.33332999999999
$ws.Range("N38").Value = -784
$ws.Range("H45").Value = 1441.8572
$ws.Range("I45").Value = 1265
$ws.Range("J45").Value = 1677.6666
$ws.Range("K45").Value = 3795
$ws.Range("L45").Value = 5032.9998
$ws.Range("M45").Value = -3263
$ws.Range("N45").Value = -6096.9998
$ws.Range("H103").Value = 413.6
$ws.Range("I103").Value = 180
$ws.Range("J103").Value = 764
$ws.Range("K103").Value = 540
$ws.Range("L103").Value = 2292
$ws.Range("M103").Value = 339
$ws.Range("N103").Value = -4050
$ws.Range("H113").Value = 1332.3334
$ws.Range("I113").Value = 600
$ws.Range("J113").Value = 1423.875
$ws.Range("K113").Value = 1800
$ws.Range("L113").Value = 4271.625
$ws.Range("M113").Value = 370
$ws.Range("N113").Value = -8611.625
$ws.Range("H129").Value = 3052.111
$ws.Range("I129").Value = 2196
$ws.Range("J129").Value = 4122.25
$ws.Range("K129").Value = 6588
$ws.Range("L129").Value = 12366.75
$ws.Range("M129").Value = -1588
$ws.Range("N129").Value = -22366.75
$ws.Range("H131").Value = 4632.212
$ws.Range("I131").Value = 825.8
$ws.Range("J131").Value = 5311.9287
$ws.Range("K131").Value = 2477.4
$ws.Range("L131").Value = 15935.7861
$ws.Range("M131").Value = 2562.6
$ws.Range("N131").Value = -26015.7861
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2090.4167
$ws.Range("I113").Value = 1910
$ws.Range("J113").Value = 2992.5
$ws.Range("K113").Value = 1910
$ws.Range("L113").Value = 2992.5
$ws.Range("M113").Value = 260
$ws.Range("N113").Value = -7332.5
$ws.Range("H123").Value = 56956.5
$ws.Range("J123").Value = 56956.5
$ws.Range("L123").Value = 56956.5
$ws.Range("N123").Value = -61856.5
$ws.Range("H132").Value = 3566.275
$ws.Range("I132").Value = 2755.7144
$ws.Range("J132").Value = 5457.5835
$ws.Range("K132").Value = 8267.143199999999
$ws.Range("L132").Value = 16372.7505
$ws.Range("M132").Value = -5737.143199999999
$ws.Range("N132").Value = -21432.7505
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2293.5715
$ws.Range("I7").Value = 2398
$ws.Range("K7").Value = 2398
$ws.Range("M7").Value = -2286
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H61").Value = 9739.416999999999
$ws.Range("I61").Value = 6565.8096
$ws.Range("K61").Value = 6565.8096
$ws.Range("M61").Value = -6363.8096
$ws.Range("H93").Value = 1621.7
$ws.Range("I93").Value = 1551.125
$ws.Range("J93").Value = 1904
$ws.Range("K93").Value = 1551.125
$ws.Range("L93").Value = 1904
$ws.Range("M93").Value = -303.125
$ws.Range("N93").Value = -4400
$ws.Range("H113").Value = 9739.416999999999
$ws.Range("I113").Value = 6565.8096
$ws.Range("K113").Value = 6565.8096
$ws.Range("M113").Value = -4395.8096
$ws.Range("H126").Value = 2293.5715
$ws.Range("I126").Value = 2398
$ws.Range("K126").Value = 7194
$ws.Range("M126").Value = -4724
$ws.Range("H136").Value = 2304.8333
$ws.Range("I136").Value = 2153.2856
$ws.Range("J136").Value = 3365.6667
$ws.Range("K136").Value = 6459.8568
$ws.Range("L136").Value = 10097.0001
$ws.Range("M136").Value = -3909.8568
$ws.Range("N136").Value = -15197.0001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 816.61536
$ws.Range("I107").Value = 827.4
$ws.Range("K107").Value = 2482.2
$ws.Range("M107").Value = -562.1999999999998
$ws.Range("H113").Value = 393.6
$ws.Range("I113").Value = 262.42856
$ws.Range("J113").Value = 560.5454999999999
$ws.Range("K113").Value = 787.28568
$ws.Range("L113").Value = 1681.6365
$ws.Range("M113").Value = 1382.71432
$ws.Range("N113").Value = -6021.6365
